$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header for column J (Q8)
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "Q8"

# New values added to row 4 (2022-07-01 00:00:00_diff)
$ws.Range("G4").Value = 0.05440228102773284
$ws.Range("H4").Value = -0.4114840973299134
$ws.Range("I4").Value = -0.01899261237092109
$ws.Range("J4").Value = -0.06854498788710228

# New values added to row 8 (2023-07-01 00:00:00_diff)
$ws.Range("G8").Value = 0.1348525673985845
$ws.Range("H8").Value = 0.01017279677448329
$ws.Range("I8").Value = -0.03366005793130281
